# Update NATMI LR-pair data: add full ECs x {ECs,FAPs,MuSCs,Resolving-Mac}
# target-cluster combinations (rows 2-17) with refreshed TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Icam4"
$ws.Cells.Item(2, 3).Value = "Itgal"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.389221
$ws.Cells.Item(2, 8).Value = 4.167663
$ws.Cells.Item(2, 9).Value = 0.2910270461264192
$ws.Cells.Item(2, 10).Value = 0.2910270461264192
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.1050723333333333
$ws.Cells.Item(2, 14).Value = 0.315217
$ws.Cells.Item(2, 15).Value = 0.006764366142259137
$ws.Cells.Item(2, 16).Value = 0.006764366142259137
$ws.Cells.Item(2, 17).Value = 0.1459686919856667
$ws.Cells.Item(2, 18).Value = 1.313718227871
$ws.Cells.Item(2, 19).Value = 0.001968613497299238
$ws.Cells.Item(2, 20).Value = 0.001968613497299238

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Icam4"
$ws.Cells.Item(3, 3).Value = "Itgal"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.389221
$ws.Cells.Item(3, 8).Value = 4.167663
$ws.Cells.Item(3, 9).Value = 0.2910270461264192
$ws.Cells.Item(3, 10).Value = 0.2910270461264192
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.9976189999999999
$ws.Cells.Item(3, 14).Value = 2.992857
$ws.Cells.Item(3, 15).Value = 0.06422490081252995
$ws.Cells.Item(3, 16).Value = 0.06422490081252995
$ws.Cells.Item(3, 17).Value = 1.385913264799
$ws.Cells.Item(3, 18).Value = 12.473219383191
$ws.Cells.Item(3, 19).Value = 0.01869118317123285
$ws.Cells.Item(3, 20).Value = 0.01869118317123285

# Row 4: ECs -> MuSCs
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Icam4"
$ws.Cells.Item(4, 3).Value = "Itgal"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.389221
$ws.Cells.Item(4, 8).Value = 4.167663
$ws.Cells.Item(4, 9).Value = 0.2910270461264192
$ws.Cells.Item(4, 10).Value = 0.2910270461264192
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.173444
$ws.Cells.Item(4, 14).Value = 0.520332
$ws.Cells.Item(4, 15).Value = 0.01116600996625811
$ws.Cells.Item(4, 16).Value = 0.01116600996625811
$ws.Cells.Item(4, 17).Value = 0.240952047124
$ws.Cells.Item(4, 18).Value = 2.168568424116
$ws.Cells.Item(4, 19).Value = 0.003249610897498254
$ws.Cells.Item(4, 20).Value = 0.003249610897498254

# Row 5: ECs -> Resolving-Mac
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Icam4"
$ws.Cells.Item(5, 3).Value = "Itgal"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.389221
$ws.Cells.Item(5, 8).Value = 4.167663
$ws.Cells.Item(5, 9).Value = 0.2910270461264192
$ws.Cells.Item(5, 10).Value = 0.2910270461264192
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 14.25707666666667
$ws.Cells.Item(5, 14).Value = 42.77123
$ws.Cells.Item(5, 15).Value = 0.9178447230789528
$ws.Cells.Item(5, 16).Value = 0.9178447230789528
$ws.Cells.Item(5, 17).Value = 19.80623030394334
$ws.Cells.Item(5, 18).Value = 178.25607273549
$ws.Cells.Item(5, 19).Value = 0.2671176385603888
$ws.Cells.Item(5, 20).Value = 0.2671176385603888

# Row 6: FAPs -> ECs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Icam4"
$ws.Cells.Item(6, 3).Value = "Itgal"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.202155333333334
$ws.Cells.Item(6, 8).Value = 6.606466000000001
$ws.Cells.Item(6, 9).Value = 0.461328155686921
$ws.Cells.Item(6, 10).Value = 0.4613281556869209
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.1050723333333333
$ws.Cells.Item(6, 14).Value = 0.315217
$ws.Cells.Item(6, 15).Value = 0.006764366142259137
$ws.Cells.Item(6, 16).Value = 0.006764366142259137
$ws.Cells.Item(6, 17).Value = 0.2313855992357778
$ws.Cells.Item(6, 18).Value = 2.082470393122001
$ws.Cells.Item(6, 19).Value = 0.00312059255679946
$ws.Cells.Item(6, 20).Value = 0.00312059255679946

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Icam4"
$ws.Cells.Item(7, 3).Value = "Itgal"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.202155333333334
$ws.Cells.Item(7, 8).Value = 6.606466000000001
$ws.Cells.Item(7, 9).Value = 0.461328155686921
$ws.Cells.Item(7, 10).Value = 0.4613281556869209
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.9976189999999999
$ws.Cells.Item(7, 14).Value = 2.992857
$ws.Cells.Item(7, 15).Value = 0.06422490081252995
$ws.Cells.Item(7, 16).Value = 0.06422490081252995
$ws.Cells.Item(7, 17).Value = 2.196912001484667
$ws.Cells.Item(7, 18).Value = 19.772208013362
$ws.Cells.Item(7, 19).Value = 0.02962875504101987
$ws.Cells.Item(7, 20).Value = 0.02962875504101987

# Row 8: FAPs -> MuSCs
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Icam4"
$ws.Cells.Item(8, 3).Value = "Itgal"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 2.202155333333334
$ws.Cells.Item(8, 8).Value = 6.606466000000001
$ws.Cells.Item(8, 9).Value = 0.461328155686921
$ws.Cells.Item(8, 10).Value = 0.4613281556869209
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.173444
$ws.Cells.Item(8, 14).Value = 0.520332
$ws.Cells.Item(8, 15).Value = 0.01116600996625811
$ws.Cells.Item(8, 16).Value = 0.01116600996625811
$ws.Cells.Item(8, 17).Value = 0.3819506296346667
$ws.Cells.Item(8, 18).Value = 3.437555666712001
$ws.Cells.Item(8, 19).Value = 0.00515119478411563
$ws.Cells.Item(8, 20).Value = 0.00515119478411563

# Row 9: FAPs -> Resolving-Mac
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Icam4"
$ws.Cells.Item(9, 3).Value = "Itgal"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 2.202155333333334
$ws.Cells.Item(9, 8).Value = 6.606466000000001
$ws.Cells.Item(9, 9).Value = 0.461328155686921
$ws.Cells.Item(9, 10).Value = 0.4613281556869209
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 14.25707666666667
$ws.Cells.Item(9, 14).Value = 42.77123
$ws.Cells.Item(9, 15).Value = 0.9178447230789528
$ws.Cells.Item(9, 16).Value = 0.9178447230789528
$ws.Cells.Item(9, 17).Value = 31.39629741924223
$ws.Cells.Item(9, 18).Value = 282.5666767731801
$ws.Cells.Item(9, 19).Value = 0.423427613304986
$ws.Cells.Item(9, 20).Value = 0.423427613304986

# Row 10: MuSCs -> ECs
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Icam4"
$ws.Cells.Item(10, 3).Value = "Itgal"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.1506176666666667
$ws.Cells.Item(10, 8).Value = 0.451853
$ws.Cells.Item(10, 9).Value = 0.03155280162368235
$ws.Cells.Item(10, 10).Value = 0.03155280162368235
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.1050723333333333
$ws.Cells.Item(10, 14).Value = 0.315217
$ws.Cells.Item(10, 15).Value = 0.006764366142259137
$ws.Cells.Item(10, 16).Value = 0.006764366142259137
$ws.Cells.Item(10, 17).Value = 0.01582574967788889
$ws.Cells.Item(10, 18).Value = 0.142431747101
$ws.Cells.Item(10, 19).Value = 0.000213434702996656
$ws.Cells.Item(10, 20).Value = 0.000213434702996656

# Row 11: MuSCs -> FAPs
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Icam4"
$ws.Cells.Item(11, 3).Value = "Itgal"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.1506176666666667
$ws.Cells.Item(11, 8).Value = 0.451853
$ws.Cells.Item(11, 9).Value = 0.03155280162368235
$ws.Cells.Item(11, 10).Value = 0.03155280162368235
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.9976189999999999
$ws.Cells.Item(11, 14).Value = 2.992857
$ws.Cells.Item(11, 15).Value = 0.06422490081252995
$ws.Cells.Item(11, 16).Value = 0.06422490081252995
$ws.Cells.Item(11, 17).Value = 0.1502590460023333
$ws.Cells.Item(11, 18).Value = 1.352331414021
$ws.Cells.Item(11, 19).Value = 0.002026475554638433
$ws.Cells.Item(11, 20).Value = 0.002026475554638433

# Row 12: MuSCs -> MuSCs
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Icam4"
$ws.Cells.Item(12, 3).Value = "Itgal"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.1506176666666667
$ws.Cells.Item(12, 8).Value = 0.451853
$ws.Cells.Item(12, 9).Value = 0.03155280162368235
$ws.Cells.Item(12, 10).Value = 0.03155280162368235
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.173444
$ws.Cells.Item(12, 14).Value = 0.520332
$ws.Cells.Item(12, 15).Value = 0.01116600996625811
$ws.Cells.Item(12, 16).Value = 0.01116600996625811
$ws.Cells.Item(12, 17).Value = 0.02612373057733334
$ws.Cells.Item(12, 18).Value = 0.235113575196
$ws.Cells.Item(12, 19).Value = 0.0003523188973934021
$ws.Cells.Item(12, 20).Value = 0.0003523188973934021

# Row 13: MuSCs -> Resolving-Mac
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Icam4"
$ws.Cells.Item(13, 3).Value = "Itgal"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.1506176666666667
$ws.Cells.Item(13, 8).Value = 0.451853
$ws.Cells.Item(13, 9).Value = 0.03155280162368235
$ws.Cells.Item(13, 10).Value = 0.03155280162368235
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 14.25707666666667
$ws.Cells.Item(13, 14).Value = 42.77123
$ws.Cells.Item(13, 15).Value = 0.9178447230789528
$ws.Cells.Item(13, 16).Value = 0.9178447230789528
$ws.Cells.Item(13, 17).Value = 2.147367621021111
$ws.Cells.Item(13, 18).Value = 19.32630858919
$ws.Cells.Item(13, 19).Value = 0.02896057246865386
$ws.Cells.Item(13, 20).Value = 0.02896057246865386

# Row 14: Resolving-Mac -> ECs
$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "Icam4"
$ws.Cells.Item(14, 3).Value = "Itgal"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 1.031517666666667
$ws.Cells.Item(14, 8).Value = 3.094553
$ws.Cells.Item(14, 9).Value = 0.2160919965629775
$ws.Cells.Item(14, 10).Value = 0.2160919965629775
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.1050723333333333
$ws.Cells.Item(14, 14).Value = 0.315217
$ws.Cells.Item(14, 15).Value = 0.006764366142259137
$ws.Cells.Item(14, 16).Value = 0.006764366142259137
$ws.Cells.Item(14, 17).Value = 0.1083839681112222
$ws.Cells.Item(14, 18).Value = 0.9754557130010001
$ws.Cells.Item(14, 19).Value = 0.001461725385163783
$ws.Cells.Item(14, 20).Value = 0.001461725385163783

# Row 15: Resolving-Mac -> FAPs
$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "Icam4"
$ws.Cells.Item(15, 3).Value = "Itgal"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 1.031517666666667
$ws.Cells.Item(15, 8).Value = 3.094553
$ws.Cells.Item(15, 9).Value = 0.2160919965629775
$ws.Cells.Item(15, 10).Value = 0.2160919965629775
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 0.9976189999999999
$ws.Cells.Item(15, 14).Value = 2.992857
$ws.Cells.Item(15, 15).Value = 0.06422490081252995
$ws.Cells.Item(15, 16).Value = 0.06422490081252995
$ws.Cells.Item(15, 17).Value = 1.029061623102333
$ws.Cells.Item(15, 18).Value = 9.261554607921
$ws.Cells.Item(15, 19).Value = 0.01387848704563879
$ws.Cells.Item(15, 20).Value = 0.01387848704563879

# Row 16: Resolving-Mac -> MuSCs
$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "Icam4"
$ws.Cells.Item(16, 3).Value = "Itgal"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 1.031517666666667
$ws.Cells.Item(16, 8).Value = 3.094553
$ws.Cells.Item(16, 9).Value = 0.2160919965629775
$ws.Cells.Item(16, 10).Value = 0.2160919965629775
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.173444
$ws.Cells.Item(16, 14).Value = 0.520332
$ws.Cells.Item(16, 15).Value = 0.01116600996625811
$ws.Cells.Item(16, 16).Value = 0.01116600996625811
$ws.Cells.Item(16, 17).Value = 0.1789105501773333
$ws.Cells.Item(16, 18).Value = 1.610194951596
$ws.Cells.Item(16, 19).Value = 0.00241288538725082
$ws.Cells.Item(16, 20).Value = 0.00241288538725082

# Row 17: Resolving-Mac -> Resolving-Mac
$ws.Cells.Item(17, 1).Value = "Resolving-Mac"
$ws.Cells.Item(17, 2).Value = "Icam4"
$ws.Cells.Item(17, 3).Value = "Itgal"
$ws.Cells.Item(17, 4).Value = "Resolving-Mac"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 1.031517666666667
$ws.Cells.Item(17, 8).Value = 3.094553
$ws.Cells.Item(17, 9).Value = 0.2160919965629775
$ws.Cells.Item(17, 10).Value = 0.2160919965629775
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 14.25707666666667
$ws.Cells.Item(17, 14).Value = 42.77123
$ws.Cells.Item(17, 15).Value = 0.9178447230789528
$ws.Cells.Item(17, 16).Value = 0.9178447230789528
$ws.Cells.Item(17, 17).Value = 14.70642645668778
$ws.Cells.Item(17, 18).Value = 132.35783811019
$ws.Cells.Item(17, 19).Value = 0.1983388987449241
$ws.Cells.Item(17, 20).Value = 0.1983388987449241

Write-Output "Updated rows 2-17 with new TPM values"
